$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.189.18'
$ws.Range('E2').Value = '  -4.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.659.32'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.34'
$ws.Range('E5').Value = '  -2.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5129'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('E8').Value = '  -3.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06440'
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.98'
$ws.Range('E10').Value = '  -3.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07813'
$ws.Range('E11').Value = '  +2.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.661.47'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.295'
$ws.Range('E13').Value = '  -4.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.887.42'
$ws.Range('E14').Value = '  -2.78%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5545'
$ws.Range('E15').Value = '  -4.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8059'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.29'
$ws.Range('E17').Value = '  -4.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.223.43'
$ws.Range('E18').Value = '  -4.03%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.78'
$ws.Range('E20').Value = '  -2.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.426'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.08'
$ws.Range('E22').Value = '  -3.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.049'
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.22'
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('E26').Value = '  +2.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1170'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.990'
$ws.Range('E28').Value = '  -3.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.79'
$ws.Range('E29').Value = '  -2.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05209'
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.252'
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.355'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.229'
$ws.Range('E33').Value = '  -5.32%  '
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.375'
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9307'
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.174.23'
$ws.Range('E38').Value = '  +12.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5706'
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01595'
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8400'
$ws.Range('E42').Value = '  -0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.681'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.70'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.797.75'
$ws.Range('E46').Value = '  +4.60%  '
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.97'
$ws.Range('E48').Value = '  -3.30%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.891'
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05063'
$ws.Range('E51').Value = '  -3.22%  '
